# MoonLake Immunotherapeutics.xlsx - "Add files via upload"
#
# Updates the hidradenitis suppurativa competitive-trial commentary block
# on the "Main" sheet: refreshes the source link, swaps the PIIb trial
# reference for the M1095-HS-201 study details, adds a new headline result
# plus comparator drug names, and records both the old and new patient
# counts for the NCT03384745 study.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New comparator callouts in column N (rows 2, 5, 6)
$ws.Range("N2").Value2 = "https://conexiant.com/dermatology/articles/sonelokimab-biologics-show-benefit-in-hidradenitis-suppurativa/"
$ws.Range("N5").Value2 = "adalimumab "
$ws.Range("N6").Value2 = "Bimekizumab "

# New headline trial result under the "In Thousands" label
$ws.Range("E10").Value2 = "Primary endpoint HiSCR75 met with 29 percentage points (ppt) delta vs placebo (p=0.0002) at week 12, setting a new bar in HS"

# Replace the old NCT03384745 results link with the new M1095-HS-201 study id/link
$ws.Range("F8").Value2 = "M1095-HS-201"
$ws.Range("F9").Value2 = "https://clinicaltrials.gov/study/NCT05322473"

# New row with the (trimmed) original clinicaltrials.gov link and its patient count
$ws.Range("E12").Value2 = "https://clinicaltrials.gov/study/NCT03384745"
$ws.Range("E9").Value2 = "n = 214"
$ws.Range("E13").Value2 = "n=313"

# Column E was given an explicit (default-sized) custom width
$ws.Columns.Item(5).ColumnWidth = 8.3

# Restore the last-used selection to match the saved workbook state
[void]$ws.Range("P18").Select()
